$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting already used by the data rows above (B2:B26) down onto
# the three new rows, the same way a user filling down a column in Excel
# would do it (Copy + Paste Formats), so the new cells pick up the existing
# "Linea" number style instead of the stray right-aligned format that used
# to sit on the old trailing blank row (B27).
$ws.Range("B26").Copy() | Out-Null
$ws.Range("B27:B29").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Lineas para Canjes realizados
$ws.Range("A27").Value = "Canjes_Realizados_MIX"
$ws.Range("B27").Value = 1162816939

$ws.Range("A28").Value = "Canjes_Realizados_POS"
$ws.Range("B28").Value = 1145642605

$ws.Range("A29").Value = "Canjes_Realizados_PRE"
$ws.Range("B29").Value = 1162676705

# Leave the view scrolled/selected the way it was left after entering the
# new rows.
$ws.Range("C27").Select()
